# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 22
    $ws.Range("F4").Value = 97
    $ws.Range("F6").Value = 47
    $ws.Range("F7").Value = 2690

    if ($sheetName -eq "展览") {
        $ws.Range("F9").Value = 273
        $ws.Range("F10").Value = 124
        $ws.Range("F11").Value = 10145
        $ws.Range("F14").Value = 11
        $ws.Range("F15").Value = 636
        $ws.Range("F16").Value = 11791
        $ws.Range("F17").Value = 12176
        $ws.Range("F18").Value = 26
        $ws.Range("F22").Value = 68
    }
    else {
        $ws.Range("F10").Value = 273
        $ws.Range("F11").Value = 124
        $ws.Range("F12").Value = 10145
        $ws.Range("F15").Value = 11
        $ws.Range("F16").Value = 636
        $ws.Range("F17").Value = 11791
        $ws.Range("F18").Value = 12176
        $ws.Range("F19").Value = 26
        $ws.Range("F23").Value = 68
    }
}

Write-Host "Done updating F-column counts."
